$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 87
$ws.Range("H87").Value2 = 15796.78
$ws.Range("J87").Value2 = 15796.78
$ws.Range("L87").Value2 = 15796.78
$ws.Range("N87").Value2 = -18292.78
# Row 90
$ws.Range("H90").Value2 = 15796.78
$ws.Range("J90").Value2 = 15796.78
$ws.Range("L90").Value2 = 47390.34
$ws.Range("N90").Value2 = -59870.34
# Row 138
$ws.Range("H138").Value2 = 3152.9
$ws.Range("I138").Value2 = 1922.9166
$ws.Range("J138").Value2 = 3844.7656
$ws.Range("K138").Value2 = 5768.7498
$ws.Range("L138").Value2 = 11534.2968
$ws.Range("M138").Value2 = -628.7497999999996
$ws.Range("N138").Value2 = -21814.2968

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value2 = 1165.4546
$ws.Range("I2").Value2 = 1030
$ws.Range("J2").Value2 = 1775
$ws.Range("K2").Value2 = 1030
$ws.Range("L2").Value2 = 1775
$ws.Range("M2").Value2 = -917
$ws.Range("N2").Value2 = -2001
# Row 110
$ws.Range("H110").Value2 = 167066.67
$ws.Range("I110").Value2 = 250100
$ws.Range("J110").Value2 = 1000
$ws.Range("K110").Value2 = 250100
$ws.Range("L110").Value2 = 1000
$ws.Range("M110").Value2 = -248055
$ws.Range("N110").Value2 = -5090
# Row 116
$ws.Range("H116").Value2 = 1165.4546
$ws.Range("I116").Value2 = 1030
$ws.Range("J116").Value2 = 1775
$ws.Range("K116").Value2 = 1030
$ws.Range("L116").Value2 = 1775
$ws.Range("M116").Value2 = 1264
$ws.Range("N116").Value2 = -6363
# Row 122
$ws.Range("H122").Value2 = 1985.5
$ws.Range("I122").Value2 = 1713.28
$ws.Range("J122").Value2 = 2957.7144
$ws.Range("K122").Value2 = 5139.84
$ws.Range("L122").Value2 = 8873.143199999999
$ws.Range("M122").Value2 = -2689.84
$ws.Range("N122").Value2 = -13773.1432

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value2 = 1165.4546
$ws.Range("I3").Value2 = 1030
$ws.Range("J3").Value2 = 1775
$ws.Range("K3").Value2 = 1030
$ws.Range("L3").Value2 = 1775
$ws.Range("M3").Value2 = -916
$ws.Range("N3").Value2 = -2003
# Row 22
$ws.Range("H22").Value2 = 267
$ws.Range("I22").Value2 = 249.625
$ws.Range("K22").Value2 = 249.625
$ws.Range("M22").Value2 = -76.625
# Row 105
$ws.Range("H105").Value2 = 2559.389
$ws.Range("I105").Value2 = 2439.2
$ws.Range("J105").Value2 = 3160.3333
$ws.Range("K105").Value2 = 2439.2
$ws.Range("L105").Value2 = 3160.3333
$ws.Range("M105").Value2 = -692.1999999999998
$ws.Range("N105").Value2 = -6654.3333
# Row 109
$ws.Range("H109").Value2 = 39245
$ws.Range("J109").Value2 = 39245
$ws.Range("L109").Value2 = 39245
$ws.Range("N109").Value2 = -42019
# Row 115
$ws.Range("H115").Value2 = 35684
$ws.Range("J115").Value2 = 35684
$ws.Range("L115").Value2 = 35684
$ws.Range("N115").Value2 = -38818

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value2 = 2011.5714
$ws.Range("I16").Value2 = 2063.5
$ws.Range("J16").Value2 = 1700
$ws.Range("K16").Value2 = 2063.5
$ws.Range("L16").Value2 = 1700
$ws.Range("M16").Value2 = -1776.5
$ws.Range("N16").Value2 = -2274
# Row 99
$ws.Range("H99").Value2 = 35767.332
$ws.Range("I99").Value2 = 2302
$ws.Range("J99").Value2 = 52500
$ws.Range("K99").Value2 = 2302
$ws.Range("L99").Value2 = 52500
$ws.Range("M99").Value2 = -804
$ws.Range("N99").Value2 = -55496
# Row 113
$ws.Range("H113").Value2 = 2011.5714
$ws.Range("I113").Value2 = 2063.5
$ws.Range("J113").Value2 = 1700
$ws.Range("K113").Value2 = 2063.5
$ws.Range("L113").Value2 = 1700
$ws.Range("M113").Value2 = 106.5
$ws.Range("N113").Value2 = -6040
# Row 122
$ws.Range("H122").Value2 = 911396.9399999999
$ws.Range("I122").Value2 = 2312
$ws.Range("J122").Value2 = 1113415.8
$ws.Range("K122").Value2 = 6936
$ws.Range("L122").Value2 = 3340247.4
$ws.Range("M122").Value2 = -4486
$ws.Range("N122").Value2 = -3345147.4
# Row 126
$ws.Range("H126").Value2 = 35767.332
$ws.Range("I126").Value2 = 2302
$ws.Range("J126").Value2 = 52500
$ws.Range("K126").Value2 = 6906
$ws.Range("L126").Value2 = 157500
$ws.Range("M126").Value2 = -4436
$ws.Range("N126").Value2 = -162440
# Row 133
$ws.Range("H133").Value2 = 0
$ws.Range("J133").Value2 = 0
$ws.Range("L133").Value2 = 0
$ws.Range("N133").ClearContents()
# Row 134
$ws.Range("H134").Value2 = 4333.2
$ws.Range("I134").Value2 = 4333.2
$ws.Range("J134").Value2 = 0
$ws.Range("K134").Value2 = 12999.6
$ws.Range("L134").Value2 = 0
$ws.Range("M134").Value2 = -10464.6
$ws.Range("N134").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 40
$ws.Range("H40").Value2 = 390.74075
$ws.Range("I40").Value2 = 78.57143000000001
$ws.Range("J40").Value2 = 500
$ws.Range("K40").Value2 = 314.28572
$ws.Range("L40").Value2 = 2000
$ws.Range("M40").Value2 = -245.28572
$ws.Range("N40").Value2 = -2138
# Row 63
$ws.Range("H63").Value2 = 3429.238
$ws.Range("I63").Value2 = 2000
$ws.Range("J63").Value2 = 3500.7
$ws.Range("K63").Value2 = 6000
$ws.Range("L63").Value2 = 10502.1
$ws.Range("M63").Value2 = -5251
$ws.Range("N63").Value2 = -12000.1
# Row 64
$ws.Range("H64").Value2 = 1948.8889
$ws.Range("I64").Value2 = 2956
$ws.Range("K64").Value2 = 8868
$ws.Range("M64").Value2 = -8598
# Row 66
$ws.Range("H66").Value2 = 3429.238
$ws.Range("I66").Value2 = 2000
$ws.Range("J66").Value2 = 3500.7
$ws.Range("K66").Value2 = 18000
$ws.Range("L66").Value2 = 31506.3
$ws.Range("M66").Value2 = -14256
$ws.Range("N66").Value2 = -38994.3
# Row 67
$ws.Range("H67").Value2 = 1948.8889
$ws.Range("I67").Value2 = 2956
$ws.Range("K67").Value2 = 8868
$ws.Range("M67").Value2 = -7932
# Row 114
$ws.Range("H114").Value2 = 299.53845
$ws.Range("I114").Value2 = 298.5
$ws.Range("J114").Value2 = 300
$ws.Range("K114").Value2 = 895.5
$ws.Range("L114").Value2 = 900
$ws.Range("M114").Value2 = 2358.5
$ws.Range("N114").Value2 = -7408
# Row 129
$ws.Range("H129").Value2 = 1440.4736
$ws.Range("I129").Value2 = 810.7692
$ws.Range("K129").Value2 = 2432.3076
$ws.Range("M129").Value2 = 2567.6924
# Row 131
$ws.Range("H131").Value2 = 841.71
$ws.Range("I131").Value2 = 315
$ws.Range("K131").Value2 = 945
$ws.Range("M131").Value2 = 4095

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value2 = 688.7143
$ws.Range("I97").Value2 = 577.5
$ws.Range("J97").Value2 = 837
$ws.Range("K97").Value2 = 577.5
$ws.Range("L97").Value2 = 837
$ws.Range("M97").Value2 = -81.5
$ws.Range("N97").Value2 = -1829
# Row 138
$ws.Range("H138").Value2 = 43752.145
$ws.Range("J138").Value2 = 43752.145
$ws.Range("L138").Value2 = 43752.145
$ws.Range("N138").Value2 = -54032.145

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value2 = 1280
$ws.Range("I61").Value2 = 1490
$ws.Range("J61").Value2 = 1000
$ws.Range("K61").Value2 = 1490
$ws.Range("L61").Value2 = 1000
$ws.Range("M61").Value2 = -1288
$ws.Range("N61").Value2 = -1404
# Row 113
$ws.Range("H113").Value2 = 1280
$ws.Range("I113").Value2 = 1490
$ws.Range("J113").Value2 = 1000
$ws.Range("K113").Value2 = 1490
$ws.Range("L113").Value2 = 1000
$ws.Range("M113").Value2 = 680
$ws.Range("N113").Value2 = -5340
# Row 127
$ws.Range("H127").Value2 = 50548.125
$ws.Range("J127").Value2 = 50548.125
$ws.Range("L127").Value2 = 50548.125
$ws.Range("N127").Value2 = -60468.125
# Row 133
$ws.Range("H133").Value2 = 36347.11
$ws.Range("J133").Value2 = 36347.11
$ws.Range("L133").Value2 = 36347.11
$ws.Range("N133").Value2 = -41407.11

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value2 = 522.2222
$ws.Range("I107").Value2 = 525
$ws.Range("J107").Value2 = 520
$ws.Range("K107").Value2 = 1575
$ws.Range("L107").Value2 = 1560
$ws.Range("M107").Value2 = 345
$ws.Range("N107").Value2 = -5400
